$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2163120567375887
$ws.Range("C2").Value = 0.5059101654846335
$ws.Range("J2").Value = 0.01773049645390071
$ws.Range("P2").Value = 0.148936170212766
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("B3").Value = 0.009070294784580499
$ws.Range("C3").Value = 0.03854875283446712
$ws.Range("J3").Value = 0.0272108843537415
$ws.Range("P3").Value = 0.673469387755102
$ws.Range("S3").Value = 0.2517006802721088
$ws.Range("J4").Value = 0.1047619047619048
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2285714285714286
$ws.Range("O5").Value = 0.125
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.375
$ws.Range("B6").Value = 0.07042253521126761
$ws.Range("D6").Value = 0.01584507042253521
$ws.Range("F6").Value = 0.06161971830985916
$ws.Range("J6").Value = 0.2394366197183098
$ws.Range("O6").Value = 0.02992957746478873
$ws.Range("Q6").Value = 0.1302816901408451
$ws.Range("R6").Value = 0.07042253521126761
$ws.Range("S6").Value = 0.3820422535211268
$ws.Range("B7").Value = 0.08637236084452975
$ws.Range("D7").Value = 0.009596928982725527
$ws.Range("F7").Value = 0.04798464491362764
$ws.Range("J7").Value = 0.1343570057581574
$ws.Range("O7").Value = 0.02303262955854127
$ws.Range("Q7").Value = 0.163147792706334
$ws.Range("R7").Value = 0.07293666026871401
$ws.Range("S7").Value = 0.4625719769673705
$ws.Range("B8").Value = 0.102263856362217
$ws.Range("D8").Value = 0.01327088212334114
$ws.Range("E8").Value = 0.00078064012490242
$ws.Range("F8").Value = 0.05308352849336456
$ws.Range("J8").Value = 0.107728337236534
$ws.Range("O8").Value = 0.02263856362217018
$ws.Range("Q8").Value = 0.1608118657298985
$ws.Range("R8").Value = 0.1046057767369243
$ws.Range("S8").Value = 0.4348165495706479
$ws.Range("B9").Value = 0.08624229979466119
$ws.Range("D9").Value = 0.01642710472279261
$ws.Range("E9").Value = 0.004106776180698152
$ws.Range("F9").Value = 0.06365503080082136
$ws.Range("J9").Value = 0.1232032854209446
$ws.Range("O9").Value = 0.01642710472279261
$ws.Range("Q9").Value = 0.1601642710472279
$ws.Range("R9").Value = 0.06776180698151951
$ws.Range("S9").Value = 0.4620123203285421
$ws.Range("B10").Value = 0.1112365894974591
$ws.Range("D10").Value = 0.01976284584980237
$ws.Range("E10").Value = 0.001976284584980237
$ws.Range("F10").Value = 0.06719367588932806
$ws.Range("J10").Value = 0.1194240542066629
$ws.Range("O10").Value = 0.01496329757199322
$ws.Range("Q10").Value = 0.2162619988706945
$ws.Range("R10").Value = 0.07340485601355166
$ws.Range("S10").Value = 0.3757763975155279
$ws.Range("F11").Value = 0.001129943502824859
$ws.Range("G11").Value = 0.1570621468926554
$ws.Range("J11").Value = 0.09830508474576272
$ws.Range("K11").Value = 0.2203389830508475
$ws.Range("L11").Value = 0.5062146892655367
$ws.Range("S11").Value = 0.01694915254237288
$ws.Range("G12").Value = 0.7100213219616205
$ws.Range("J12").Value = 0.2260127931769723
$ws.Range("K12").Value = 0.006396588486140725
$ws.Range("L12").Value = 0.02132196162046908
$ws.Range("S12").Value = 0.03624733475479744
$ws.Range("F13").Value = 0.008064516129032258
$ws.Range("G13").Value = 0.5806451612903226
$ws.Range("J13").Value = 0.3387096774193548
$ws.Range("S13").Value = 0.07258064516129033
$ws.Range("F15").Value = 0.01379310344827586
$ws.Range("H15").Value = 0.1396551724137931
$ws.Range("I15").Value = 0.04655172413793104
$ws.Range("J15").Value = 0.3724137931034483
$ws.Range("K15").Value = 0.06379310344827586
$ws.Range("M15").Value = 0.008620689655172414
$ws.Range("O15").Value = 0.05517241379310345
$ws.Range("S15").Value = 0.3
$ws.Range("F16").Value = 0.01452282157676349
$ws.Range("H16").Value = 0.1742738589211618
$ws.Range("I16").Value = 0.08298755186721991
$ws.Range("J16").Value = 0.3900414937759336
$ws.Range("K16").Value = 0.1058091286307054
$ws.Range("M16").Value = 0.02489626556016597
$ws.Range("O16").Value = 0.03526970954356846
$ws.Range("S16").Value = 0.1721991701244813
$ws.Range("F17").Value = 0.01327800829875519
$ws.Range("H17").Value = 0.1634854771784232
$ws.Range("I17").Value = 0.08298755186721991
$ws.Range("J17").Value = 0.4348547717842324
$ws.Range("K17").Value = 0.0970954356846473
$ws.Range("M17").Value = 0.01659751037344398
$ws.Range("O17").Value = 0.05477178423236514
$ws.Range("S17").Value = 0.1369294605809129
$ws.Range("F18").Value = 0.02390438247011952
$ws.Range("H18").Value = 0.1474103585657371
$ws.Range("I18").Value = 0.08565737051792828
$ws.Range("J18").Value = 0.4183266932270917
$ws.Range("K18").Value = 0.1175298804780877
$ws.Range("M18").Value = 0.009960159362549801
$ws.Range("O18").Value = 0.05179282868525897
$ws.Range("S18").Value = 0.1454183266932271
$ws.Range("F19").Value = 0.01623886851754846
$ws.Range("H19").Value = 0.2239392352016763
$ws.Range("I19").Value = 0.07333682556312206
$ws.Range("J19").Value = 0.3525405971712939
$ws.Range("K19").Value = 0.1089575694080671
$ws.Range("M19").Value = 0.02226296490309062
$ws.Range("N19").Value = 0.000785751702462022
$ws.Range("O19").Value = 0.06338397066526978
$ws.Range("S19").Value = 0.1385542168674699
